# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the "Repayment
# schedule" sheet, shifting the old N ("Late"), O (blank heading) and
# P ("Outstanding") columns one place to the right, then leave the
# "Repayment schedule" sheet active with the selection on R6 (matching the
# post-insert selection recorded by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new column at N - everything from N onward shifts right by one.
$ws.Columns.Item(14).Insert()

# The freshly inserted column inherits the width Excel would give a column
# inserted to the left of "In Advance" (M) - 11 characters, not a bestFit
# width.
$ws.Columns.Item(14).ColumnWidth = 11 - 0.8333333333333333

# Final selection left on the sheet after the insert.
$ws.Range("R6").Select() | Out-Null
